$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new range as Text first so numeric-looking strings
# (e.g. "26", "136.84") are stored as text, matching the rest of the sheet.
$ws.Range("A13:K23").NumberFormat = "@"

# Row 13
$ws.Range("A13").Value = " Dubai (DSC)"
$ws.Range("B13").Value = " October 04 2020"
$ws.Range("C13").Value = "Super Kings won by 10 wickets (with 14 balls remaining)"
$ws.Range("D13").Value = "Kings XI Punjab"
$ws.Range("E13").Value = "Chennai Super Kings"
$ws.Range("F13").Value = "Mayank Agarwal "
$ws.Range("G13").Value = "26"
$ws.Range("H13").Value = "19"
$ws.Range("I13").Value = "3"
$ws.Range("J13").Value = "0"
$ws.Range("K13").Value = "136.84"

# Row 14
$ws.Range("A14").Value = " Abu Dhabi"
$ws.Range("B14").Value = " November 01 2020"
$ws.Range("C14").Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Range("D14").Value = "Kings XI Punjab"
$ws.Range("E14").Value = "Chennai Super Kings"
$ws.Range("F14").Value = "Mayank Agarwal "
$ws.Range("G14").Value = "26"
$ws.Range("H14").Value = "15"
$ws.Range("I14").Value = "5"
$ws.Range("J14").Value = "0"
$ws.Range("K14").Value = "173.33"

# Row 15
$ws.Range("A15").Value = " Abu Dhabi"
$ws.Range("B15").Value = " October 01 2020"
$ws.Range("C15").Value = "Mumbai won by 48 runs"
$ws.Range("D15").Value = "Kings XI Punjab"
$ws.Range("E15").Value = "Mumbai Indians"
$ws.Range("F15").Value = "Mayank Agarwal "
$ws.Range("G15").Value = "25"
$ws.Range("H15").Value = "18"
$ws.Range("I15").Value = "3"
$ws.Range("J15").Value = "0"
$ws.Range("K15").Value = "138.88"

# Row 16
$ws.Range("A16").Value = " Dubai (DSC)"
$ws.Range("B16").Value = " September 24 2020"
$ws.Range("C16").Value = "Kings XI won by 97 runs"
$ws.Range("D16").Value = "Kings XI Punjab"
$ws.Range("E16").Value = "Royal Challengers Bangalore"
$ws.Range("F16").Value = "Mayank Agarwal "
$ws.Range("G16").Value = "26"
$ws.Range("H16").Value = "20"
$ws.Range("I16").Value = "4"
$ws.Range("J16").Value = "0"
$ws.Range("K16").Value = "130.00"

# Row 17
$ws.Range("A17").Value = " Sharjah"
$ws.Range("B17").Value = " October 15 2020"
$ws.Range("C17").Value = "Kings XI won by 8 wickets"
$ws.Range("D17").Value = "Kings XI Punjab"
$ws.Range("E17").Value = "Royal Challengers Bangalore"
$ws.Range("F17").Value = "Mayank Agarwal "
$ws.Range("G17").Value = "45"
$ws.Range("H17").Value = "25"
$ws.Range("I17").Value = "4"
$ws.Range("J17").Value = "3"
$ws.Range("K17").Value = "180.00"

# Row 18
$ws.Range("A18").Value = " Dubai (DSC)"
$ws.Range("B18").Value = " October 20 2020"
$ws.Range("C18").Value = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Range("D18").Value = "Kings XI Punjab"
$ws.Range("E18").Value = "Delhi Capitals"
$ws.Range("F18").Value = "Mayank Agarwal "
$ws.Range("G18").Value = "5"
$ws.Range("H18").Value = "9"
$ws.Range("I18").Value = "0"
$ws.Range("J18").Value = "0"
$ws.Range("K18").Value = "55.55"

# Row 19
$ws.Range("A19").Value = " Abu Dhabi"
$ws.Range("B19").Value = " October 10 2020"
$ws.Range("C19").Value = "KKR won by 2 runs"
$ws.Range("D19").Value = "Kings XI Punjab"
$ws.Range("E19").Value = "Kolkata Knight Riders"
$ws.Range("F19").Value = "Mayank Agarwal "
$ws.Range("G19").Value = "56"
$ws.Range("H19").Value = "39"
$ws.Range("I19").Value = "6"
$ws.Range("J19").Value = "1"
$ws.Range("K19").Value = "143.58"

# Row 20
$ws.Range("A20").Value = " Dubai (DSC)"
$ws.Range("B20").Value = " September 20 2020"
$ws.Range("C20").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D20").Value = "Kings XI Punjab"
$ws.Range("E20").Value = "Delhi Capitals"
$ws.Range("F20").Value = "Mayank Agarwal "
$ws.Range("G20").Value = "89"
$ws.Range("H20").Value = "60"
$ws.Range("I20").Value = "7"
$ws.Range("J20").Value = "4"
$ws.Range("K20").Value = "148.33"

# Row 21
$ws.Range("A21").Value = " Sharjah"
$ws.Range("B21").Value = " September 27 2020"
$ws.Range("C21").Value = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Range("D21").Value = "Kings XI Punjab"
$ws.Range("E21").Value = "Rajasthan Royals"
$ws.Range("F21").Value = "Mayank Agarwal "
$ws.Range("G21").Value = "106"
$ws.Range("H21").Value = "50"
$ws.Range("I21").Value = "10"
$ws.Range("J21").Value = "7"
$ws.Range("K21").Value = "212.00"

# Row 22
$ws.Range("A22").Value = " Dubai (DSC)"
$ws.Range("B22").Value = " October 18 2020"
$ws.Range("C22").Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Range("D22").Value = "Kings XI Punjab"
$ws.Range("E22").Value = "Mumbai Indians"
$ws.Range("F22").Value = "Mayank Agarwal "
$ws.Range("G22").Value = "11"
$ws.Range("H22").Value = "10"
$ws.Range("I22").Value = "1"
$ws.Range("J22").Value = "0"
$ws.Range("K22").Value = "110.00"

# Row 23
$ws.Range("A23").Value = " Dubai (DSC)"
$ws.Range("B23").Value = " October 08 2020"
$ws.Range("C23").Value = "Sunrisers won by 69 runs"
$ws.Range("D23").Value = "Kings XI Punjab"
$ws.Range("E23").Value = "Sunrisers Hyderabad"
$ws.Range("F23").Value = "Mayank Agarwal "
$ws.Range("G23").Value = "9"
$ws.Range("H23").Value = "6"
$ws.Range("I23").Value = "1"
$ws.Range("J23").Value = "0"
$ws.Range("K23").Value = "150.00"

# Extend the "number stored as text" ignored-error hint to the new range,
# mirroring the original sheet-wide suppression (best effort; matches Excel's
# Range.Errors(xlNumberAsText).Ignore behaviour).
$ws.Range("A1:K23").Errors.Item(3).Ignore = $true

Write-Host "done"